# Update "evening debriefing" events on the Events sheet.
#
# Summary of changes (per commit "updated evening debriefing events"):
#  - Row 61 (e100): body text (col B) expanded with a reference to Crew
#    Rating Improvements (r4.91).
#  - Row 62: id corrected from "e0101" -> "e101" and body text (col B)
#    replaced with new "Victory Point Total" copy.
#  - Two brand new rows inserted after the (updated) e101 row:
#       e102 - Evening Debriefing - Promotions
#       e103 - Evening Debriefing - Decorations
#  - The following rows (previously e501, e502, e503, e503a, e504) are
#    pushed down by two rows but are otherwise unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 61 (e100 Evening Debriefing) ---------------------
$ws.Range("B61").Value = "<Bold>e100 Evening Debriefing</Bold> `n<LineBreak/><LineBreak/>`nAn evening debriefing is performed per <InlineUIContainer><Button Content='r4.9' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. `nClick image to continue to continue to Crew Rating Improvements per `n<InlineUIContainer><Button Content='r4.91' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   `n<LineBreak/><LineBreak/>`n                     <InlineUIContainer><Image Name='Debrief' Height='168' Width='275'></Image></InlineUIContainer>"
$ws.Rows.Item(61).RowHeight = 105

# --- Update existing row 62: e0101 -> e101, new body text -----------------
$ws.Range("A62").Value = "e101"
$ws.Range("B62").Value = "<Bold>e101 Evening Debriefing - Victory Point Total</Bold> `n<LineBreak/><LineBreak/>`nThe After Action Report `n<InlineUIContainer><Button Content='AAR' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `nis updated to reflect victory points for both yoru tank and friendly forces. `n<LineBreak/><LineBreak/>`nIf the combined victory points from both your tank and friendly forces is positive, you have won the engagement.`n<LineBreak/><LineBreak/>"
$ws.Rows.Item(62).RowHeight = 120

# --- Insert two new rows for e102 and e103 after row 62 --------------------
$ws.Rows.Item(63).Resize(2).Insert()

$ws.Range("A63").Value = "e102"
$ws.Range("B63").Value = "<Bold>e102 Evening Debriefing - Promotions</Bold> `n<InlineUIContainer><Button Content='r25.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `n<LineBreak/><LineBreak/>`nThe After Action Report `n<InlineUIContainer><Button Content='AAR' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `nis updated to reflect promotions.`n<LineBreak/><LineBreak/>"
$ws.Rows.Item(63).RowHeight = 105

$ws.Range("A64").Value = "e103"
$ws.Range("B64").Value = "<Bold>e103 Evening Debriefing - Decorations</Bold> `n<InlineUIContainer><Button Content='r26.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `n<LineBreak/><LineBreak/>`nRoll for possible decorations on the `n<InlineUIContainer><Button Content='Decorations' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Table.`nMedals received are recorded on the  per After Action Report `n<InlineUIContainer><Button Content='AAR' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> .`n<LineBreak/><LineBreak/>"
$ws.Rows.Item(64).RowHeight = 120

# (Rows.Insert() already copies formatting from the row above, so the new
# A63:B64 cells automatically pick up the same styles used throughout the
# rest of the A/B columns - no explicit restyle needed.)

# --- Update sheet view to focus on the newly edited area -------------------
$ws.Application.ActiveWindow.ScrollRow = 61
$ws.Range("B63").Select()
